$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 324-359: realigned weekly price records (3 new records inserted at the
# top of the block, pushing the previously-existing records down by 3 rows and
# the final 3 records into brand-new rows 357-359).
$rows = @(
    ,@(324, 6, "Mercado Mayorista Lo Valledor de Santiago", "Metropolitana", 44578, 13, 100112043, "Pepino ensalada", "Sin especificar", "Primera", 1260, 9000, 10000, 9397, "`$/caja 60 unidades", "Región de Arica y Parinacota", 157, 60, "Hortaliza")
    ,@(325, 6, "Mercado Mayorista Lo Valledor de Santiago", "Metropolitana", 44578, 13, 100112043, "Pepino ensalada", "Sin especificar", "Primera", 290, 11000, 11000, 11000, "`$/caja 70 unidades", "Región Metropolitana", 157, 70, "Hortaliza")
    ,@(326, 6, "Mercado Mayorista Lo Valledor de Santiago", "Metropolitana", 44578, 13, 100112043, "Pepino ensalada", "Sin especificar", "Segunda", 410, 5000, 5500, 5195, "`$/caja 100 unidades", "Región de Arica y Parinacota", 52, 100, "Hortaliza")
    ,@(327, 6, "Mercado Mayorista Lo Valledor de Santiago", "Metropolitana", 44490, 13, 100112043, "Pepino ensalada", "Sin especificar", "Primera", 4200, 6000, 7000, 6286, "`$/caja 60 unidades", "Región de Arica y Parinacota", 105, 60, "Hortaliza")
    ,@(328, 6, "Mercado Mayorista Lo Valledor de Santiago", "Metropolitana", 44490, 13, 100112043, "Pepino ensalada", "Sin especificar", "Segunda", 650, 5000, 6000, 5569, "`$/caja 100 unidades", "Región de Arica y Parinacota", 56, 100, "Hortaliza")
    ,@(329, 6, "Mercado Mayorista Lo Valledor de Santiago", "Metropolitana", 44427, 13, 100112043, "Pepino ensalada", "Sin especificar", "Primera", 400, 15000, 16000, 15425, "`$/caja 60 unidades", "Región de Arica y Parinacota", 257, 60, "Hortaliza")
    ,@(330, 6, "Mercado Mayorista Lo Valledor de Santiago", "Metropolitana", 44491, 13, 100112043, "Pepino ensalada", "Sin especificar", "Primera", 700, 5000, 6000, 5543, "`$/caja 60 unidades", "Región Metropolitana", 92, 60, "Hortaliza")
    ,@(331, 6, "Mercado Mayorista Lo Valledor de Santiago", "Metropolitana", 44293, 13, 100112043, "Pepino ensalada", "Sin especificar", "Primera", 250, 13000, 14000, 13520, "`$/caja 60 unidades", "Limache", 225, 60, "Hortaliza")
    ,@(332, 6, "Mercado Mayorista Lo Valledor de Santiago", "Metropolitana", 44293, 13, 100112043, "Pepino ensalada", "Sin especificar", "Primera", 270, 12000, 13000, 12481, "`$/caja 60 unidades", "Región de Arica y Parinacota", 208, 60, "Hortaliza")
    ,@(333, 6, "Mercado Mayorista Lo Valledor de Santiago", "Metropolitana", 44266, 13, 100112043, "Pepino ensalada", "Sin especificar", "Primera", 240, 13000, 13000, 13000, "`$/caja 60 unidades", "Región de Arica y Parinacota", 217, 60, "Hortaliza")
    ,@(334, 6, "Mercado Mayorista Lo Valledor de Santiago", "Metropolitana", 44266, 13, 100112043, "Pepino ensalada", "Sin especificar", "Primera", 220, 15000, 15000, 15000, "`$/caja 70 unidades", "Provincia de Quillota", 214, 70, "Hortaliza")
    ,@(335, 6, "Mercado Mayorista Lo Valledor de Santiago", "Metropolitana", 44533, 13, 100112043, "Pepino ensalada", "Sin especificar", "Primera", 200, 6000, 7000, 6600, "`$/caja 70 unidades", "Provincia de Huasco", 94, 70, "Hortaliza")
    ,@(336, 6, "Mercado Mayorista Lo Valledor de Santiago", "Metropolitana", 44264, 13, 100112043, "Pepino ensalada", "Sin especificar", "Primera", 270, 12000, 12000, 12000, "`$/caja 60 unidades", "Provincia de Quillota", 200, 60, "Hortaliza")
    ,@(337, 6, "Mercado Mayorista Lo Valledor de Santiago", "Metropolitana", 44264, 13, 100112043, "Pepino ensalada", "Sin especificar", "Primera", 220, 13000, 13000, 13000, "`$/caja 60 unidades", "Región de Arica y Parinacota", 217, 60, "Hortaliza")
    ,@(338, 6, "Mercado Mayorista Lo Valledor de Santiago", "Metropolitana", 44494, 13, 100112043, "Pepino ensalada", "Sin especificar", "Primera", 4200, 6000, 7000, 6286, "`$/caja 60 unidades", "Región de Arica y Parinacota", 105, 60, "Hortaliza")
    ,@(339, 6, "Mercado Mayorista Lo Valledor de Santiago", "Metropolitana", 44571, 13, 100112043, "Pepino ensalada", "Sin especificar", "Primera", 2000, 5000, 6000, 5400, "`$/caja 60 unidades", "Región de Arica y Parinacota", 90, 60, "Hortaliza")
    ,@(340, 6, "Mercado Mayorista Lo Valledor de Santiago", "Metropolitana", 44571, 13, 100112043, "Pepino ensalada", "Sin especificar", "Segunda", 400, 5000, 6000, 5425, "`$/caja 100 unidades", "Región de Arica y Parinacota", 54, 100, "Hortaliza")
    ,@(341, 6, "Mercado Mayorista Lo Valledor de Santiago", "Metropolitana", 44279, 13, 100112043, "Pepino ensalada", "Sin especificar", "Primera", 500, 10000, 11000, 10640, "`$/caja 60 unidades", "Región de Arica y Parinacota", 177, 60, "Hortaliza")
    ,@(342, 6, "Mercado Mayorista Lo Valledor de Santiago", "Metropolitana", 44277, 13, 100112043, "Pepino ensalada", "Sin especificar", "Primera", 500, 10000, 11000, 10360, "`$/caja 70 unidades", "Región de Arica y Parinacota", 148, 70, "Hortaliza")
    ,@(343, 6, "Mercado Mayorista Lo Valledor de Santiago", "Metropolitana", 44277, 13, 100112043, "Pepino ensalada", "Sin especificar", "Segunda", 400, 8000, 9000, 8575, "`$/caja 100 unidades", "Región de Arica y Parinacota", 86, 100, "Hortaliza")
    ,@(344, 6, "Mercado Mayorista Lo Valledor de Santiago", "Metropolitana", 44525, 13, 100112043, "Pepino ensalada", "Sin especificar", "Primera", 3500, 5000, 6000, 5429, "`$/caja 60 unidades", "Región de Arica y Parinacota", 90, 60, "Hortaliza")
    ,@(345, 6, "Mercado Mayorista Lo Valledor de Santiago", "Metropolitana", 44525, 13, 100112043, "Pepino ensalada", "Sin especificar", "Segunda", 800, 4000, 5000, 4438, "`$/caja 100 unidades", "Región de Arica y Parinacota", 44, 100, "Hortaliza")
    ,@(346, 6, "Mercado Mayorista Lo Valledor de Santiago", "Metropolitana", 44327, 13, 100112043, "Pepino ensalada", "Sin especificar", "Primera", 290, 8000, 9000, 8586, "`$/caja 60 unidades", "Región Metropolitana", 143, 60, "Hortaliza")
    ,@(347, 6, "Mercado Mayorista Lo Valledor de Santiago", "Metropolitana", 44354, 13, 100112043, "Pepino ensalada", "Sin especificar", "Primera", 1220, 9000, 10000, 9377, "`$/caja 60 unidades", "Región de Arica y Parinacota", 156, 60, "Hortaliza")
    ,@(348, 6, "Mercado Mayorista Lo Valledor de Santiago", "Metropolitana", 44354, 13, 100112043, "Pepino ensalada", "Sin especificar", "Segunda", 300, 7000, 8000, 7333, "`$/caja 100 unidades", "Región de Arica y Parinacota", 73, 100, "Hortaliza")
    ,@(349, 6, "Mercado Mayorista Lo Valledor de Santiago", "Metropolitana", 44503, 13, 100112043, "Pepino ensalada", "Sin especificar", "Primera", 1100, 6000, 7000, 6409, "`$/caja 60 unidades", "Región de Arica y Parinacota", 107, 60, "Hortaliza")
    ,@(350, 6, "Mercado Mayorista Lo Valledor de Santiago", "Metropolitana", 44503, 13, 100112043, "Pepino ensalada", "Sin especificar", "Primera", 200, 7000, 8000, 7400, "`$/caja 70 unidades", "Provincia de Huasco", 106, 70, "Hortaliza")
    ,@(351, 6, "Mercado Mayorista Lo Valledor de Santiago", "Metropolitana", 44462, 13, 100112043, "Pepino ensalada", "Sin especificar", "Primera", 650, 14000, 15000, 14646, "`$/caja 60 unidades", "Región de Arica y Parinacota", 244, 60, "Hortaliza")
    ,@(352, 6, "Mercado Mayorista Lo Valledor de Santiago", "Metropolitana", 44462, 13, 100112043, "Pepino ensalada", "Sin especificar", "Primera", 200, 15000, 16000, 15600, "`$/caja 70 unidades", "Provincia de Huasco", 223, 70, "Hortaliza")
    ,@(353, 6, "Mercado Mayorista Lo Valledor de Santiago", "Metropolitana", 44462, 13, 100112043, "Pepino ensalada", "Sin especificar", "Segunda", 180, 12000, 13000, 12556, "`$/caja 100 unidades", "Provincia de Huasco", 126, 100, "Hortaliza")
    ,@(354, 6, "Mercado Mayorista Lo Valledor de Santiago", "Metropolitana", 44384, 13, 100112043, "Pepino ensalada", "Sin especificar", "Primera", 500, 14000, 15000, 14600, "`$/caja 60 unidades", "Región de Arica y Parinacota", 243, 60, "Hortaliza")
    ,@(355, 6, "Mercado Mayorista Lo Valledor de Santiago", "Metropolitana", 44384, 13, 100112043, "Pepino ensalada", "Sin especificar", "Segunda", 210, 12000, 13000, 12476, "`$/caja 100 unidades", "Región de Arica y Parinacota", 125, 100, "Hortaliza")
    ,@(356, 6, "Mercado Mayorista Lo Valledor de Santiago", "Metropolitana", 44512, 13, 100112043, "Pepino ensalada", "Sin especificar", "Primera", 1800, 5000, 6000, 5444, "`$/caja 60 unidades", "Región de Arica y Parinacota", 91, 60, "Hortaliza")
    ,@(357, 6, "Mercado Mayorista Lo Valledor de Santiago", "Metropolitana", 44312, 13, 100112043, "Pepino ensalada", "Sin especificar", "Primera", 400, 8000, 9000, 8575, "`$/caja 60 unidades", "Región de Arica y Parinacota", 143, 60, "Hortaliza")
    ,@(358, 6, "Mercado Mayorista Lo Valledor de Santiago", "Metropolitana", 44312, 13, 100112043, "Pepino ensalada", "Sin especificar", "Segunda", 400, 7000, 8000, 7425, "`$/caja 100 unidades", "Región de Arica y Parinacota", 74, 100, "Hortaliza")
    ,@(359, 6, "Mercado Mayorista Lo Valledor de Santiago", "Metropolitana", 44511, 13, 100112043, "Pepino ensalada", "Sin especificar", "Primera", 1300, 5000, 6000, 5577, "`$/caja 50 unidades", "Región de Arica y Parinacota", 112, 50, "Hortaliza")
)

foreach ($row in $rows) {
    $r = $row[0]
    for ($i = 0; $i -lt 18; $i++) {
        $col = $i + 1
        $ws.Cells.Item($r, $col).Value = $row[$i + 1]
    }
}

# Column D holds dates; make sure the date number format is applied to the
# whole block (existing rows already have it, new rows 357:359 need it set).
$ws.Range("D324:D359").NumberFormat = "YYYY-MM-DD HH:MM:SS"
